# Auto-generated edit script applying cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.132.90"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "1.631.77"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "1.860.55"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "1.627.86"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "27.103.06"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "1.308.65"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.805"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "1.769.23"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.812"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +20.09%  "
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("E51").Value = "  -1.21%  "
